$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 271; this shifts rows 271-340 down to 272-341
$ws.Rows.Item(271).Insert()

# Populate the new row 271 with the new record's data
$ws.Cells.Item(271, 1).Value = 6
$ws.Cells.Item(271, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(271, 3).Value = "Metropolitana"
$ws.Cells.Item(271, 4).Value = 44551
$ws.Cells.Item(271, 5).Value = 13
$ws.Cells.Item(271, 6).Value = 100112039
$ws.Cells.Item(271, 7).Value = "Ciboulette"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 1010
$ws.Cells.Item(271, 11).Value = 700
$ws.Cells.Item(271, 12).Value = 800
$ws.Cells.Item(271, 13).Value = 756
$ws.Cells.Item(271, 14).Value = "$/docena de atados"
$ws.Cells.Item(271, 15).Value = "Región Metropolitana"
$ws.Cells.Item(271, 16).Value = 252
$ws.Cells.Item(271, 17).Value = 3
$ws.Cells.Item(271, 18).Value = "Hortaliza"
